$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells (Coin name / Link / Price / Volume columns) to match
# the latest scraped crypto data. Numeric-looking Price values are written
# with a leading apostrophe (forces text, matching original t="inlineStr"/
# shared-string cell type) and then the cell style is reset to "Normal" so
# no stray quote-prefix style / number format is left behind on the cell.

$ws.Range("D2").Value = "42.270.14"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "2.174.79"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'253.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.26%  "

$ws.Range("E6").Value = "  +0.89%  "

$ws.Range("D7").Value = "'74.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("D10").Value = "'40.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("D11").Value = "'0.0908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.101"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").Value = "2.497.92"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "'14.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").Value = "2.152.90"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("D17").Value = "'0.765"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("D18").Value = "42.195.70"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").Value = "'70.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("D21").Value = "'5.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("D22").Value = "'226.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("E23").Value = "  +5.46%  "

$ws.Range("D24").Value = "'9.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.30%  "

$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("D26").Value = "'10.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").Value = "'3.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.51%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'37.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.37%  "

$ws.Range("D31").Value = "'168.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("D32").Value = "'19.97"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = "'0.0806"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.88%  "

$ws.Range("D34").Value = "'5.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.68%  "

$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("E36").Value = "  +4.23%  "

$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("D38").Value = "'0.0334"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.06%  "

$ws.Range("D39").Value = "'11.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.06%  "

$ws.Range("E40").Value = "  -2.23%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'59.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.195"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.46%  "

$ws.Range("D43").Value = "'5.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.30%  "

$ws.Range("D44").Value = "'102.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.42%  "

$ws.Range("D45").Value = "'0.469"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +16.61%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.55%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0970"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.30%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.14%  "

$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D50").Value = "'1.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("E51").Value = "  +0.64%  "
